# V 2.0.2 se arreglo la fechar y hora de reimpresion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a value as literal text, avoiding Excel's automatic
# number/date reinterpretation (e.g. "1988-12-19" -> date serial,
# "46491184" -> number). We build the text via a formula that evaluates
# to the exact string, then copy/paste-special as values only, which
# bakes the literal text into the cell without touching its style.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# Nombre del paciente / No. Expediente Clínico
$ws.Range("A6").Value = "BARRIOS   DE LEÓN  JENNIFER  MAYONI"
$ws.Range("G6").Value = "8304/201762164"

# Fecha de nacimiento / Edad / Lugar de nacimiento
Set-TextValue $ws.Range("A9") "1988-12-19"
Set-TextValue $ws.Range("D9") "28"
$ws.Range("E9").Value = "SAN JOSÉ EL RODEO SAN MARCOS"

# Documento de identificación
$ws.Range("G11").Value = "DPI 1597776821214"

# Datos de la persona a notificar en caso de emergencia
$ws.Range("A13").Value = "SAMUEL SY"
$ws.Range("D13").Value = "ESPOSO"
$ws.Range("E13").Value = "LOTE 18 MANZ. A SECT. SAN PEDRO AYAMPUC"
Set-TextValue $ws.Range("G13") "46491184"

# Hora de la asistencia médica (reimpresión)
$ws.Range("D14").Value = "Hora: 19:33:8"
